$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-10-23 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-10-24 Thursday", 2) | Out-Null
$d.Content.Find.Execute("334×2=668", $true, $false, $false, $false, $false, $true, 1, $false, "817×9=7353", 2) | Out-Null
$d.Content.Find.Execute("905×4=3620", $true, $false, $false, $false, $false, $true, 1, $false, "978×6=5868", 2) | Out-Null
$d.Content.Find.Execute("342×8=2736", $true, $false, $false, $false, $false, $true, 1, $false, "982×6=5892", 2) | Out-Null
$d.Content.Find.Execute("226×6=1356", $true, $false, $false, $false, $false, $true, 1, $false, "775×2=1550", 2) | Out-Null
$d.Content.Find.Execute("879×9=7911", $true, $false, $false, $false, $false, $true, 1, $false, "784×5=3920", 2) | Out-Null
$d.Content.Find.Execute("699×9=6291", $true, $false, $false, $false, $false, $true, 1, $false, "392×2=784", 2) | Out-Null
$d.Content.Find.Execute("606×7=4242", $true, $false, $false, $false, $false, $true, 1, $false, "707×6=4242", 2) | Out-Null
$d.Content.Find.Execute("903×3=2709", $true, $false, $false, $false, $false, $true, 1, $false, "202×9=1818", 2) | Out-Null
$d.Content.Find.Execute("375×6=2250", $true, $false, $false, $false, $false, $true, 1, $false, "678×3=2034", 2) | Out-Null
$d.Content.Find.Execute("259×6=1554", $true, $false, $false, $false, $false, $true, 1, $false, "418×8=3344", 2) | Out-Null
$d.Content.Find.Execute("857×6=5142", $true, $false, $false, $false, $false, $true, 1, $false, "463×4=1852", 2) | Out-Null
$d.Content.Find.Execute("136×7=952", $true, $false, $false, $false, $false, $true, 1, $false, "273×3=819", 2) | Out-Null
$d.Content.Find.Execute("568×6=3408", $true, $false, $false, $false, $false, $true, 1, $false, "838×8=6704", 2) | Out-Null
$d.Content.Find.Execute("263×5=1315", $true, $false, $false, $false, $false, $true, 1, $false, "194×9=1746", 2) | Out-Null
$d.Content.Find.Execute("384×2=768", $true, $false, $false, $false, $false, $true, 1, $false, "502×8=4016", 2) | Out-Null
$d.Content.Find.Execute("473×3=1419", $true, $false, $false, $false, $false, $true, 1, $false, "518×9=4662", 2) | Out-Null
$d.Content.Find.Execute("774×4=3096", $true, $false, $false, $false, $false, $true, 1, $false, "102×8=816", 2) | Out-Null
$d.Content.Find.Execute("953×9=8577", $true, $false, $false, $false, $false, $true, 1, $false, "489×3=1467", 2) | Out-Null
$d.Content.Find.Execute("775×5=3875", $true, $false, $false, $false, $false, $true, 1, $false, "642×4=2568", 2) | Out-Null
$d.Content.Find.Execute("687×6=4122", $true, $false, $false, $false, $false, $true, 1, $false, "152×7=1064", 2) | Out-Null
$d.Content.Find.Execute("841×2=1682", $true, $false, $false, $false, $false, $true, 1, $false, "383×6=2298", 2) | Out-Null
$d.Content.Find.Execute("814×3=2442", $true, $false, $false, $false, $false, $true, 1, $false, "407×4=1628", 2) | Out-Null
$d.Content.Find.Execute("704×4=2816", $true, $false, $false, $false, $false, $true, 1, $false, "479×9=4311", 2) | Out-Null
$d.Content.Find.Execute("480×2=960", $true, $false, $false, $false, $false, $true, 1, $false, "991×7=6937", 2) | Out-Null
$d.Content.Find.Execute("612×2=1224", $true, $false, $false, $false, $false, $true, 1, $false, "659×3=1977", 2) | Out-Null
